$wb = $excel.ActiveWorkbook

# Insert a new "employment" worksheet immediately before "raw data"
# (i.e. after "partnership"), mirroring the existing "partnership" sheet's
# simple year/value layout.
$rawData = $wb.Worksheets.Item("raw data")
$newSheet = $wb.Worksheets.Add($rawData)
$newSheet.Name = "employment"
$newSheet.Activate()

# Header row
$newSheet.Range("A1").Value = "year"
$newSheet.Range("B1").Value = "employed_share"

# Draft activity alignment: constant employed share of 0.6 for 2010-2027
$years = 2010..2027
for ($i = 0; $i -lt $years.Count; $i++) {
    $row = $i + 2
    $newSheet.Cells.Item($row, 1).Value = $years[$i]
    $newSheet.Cells.Item($row, 2).Value = 0.6
}

$newSheet.Range("B2").Select()
